# Script: apply the "05-11-2023" betexplorer scrape update to the
# Denmark 2nd-division 2023-2024 sheet.
#
#  1) Four existing row pairs had their match data (columns F:V) swapped
#     between two adjacent rows (the "Indice"/date columns A:E stayed put).
#  2) Four brand-new match rows (80-83) were appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, $rowA, $rowB, $firstCol, $lastCol)

    $valsA = @()
    $valsB = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $valsA += ,$ws.Cells.Item($rowA, $c).Value()
        $valsB += ,$ws.Cells.Item($rowB, $c).Value()
    }

    $c = $firstCol
    foreach ($v in $valsB) {
        $ws.Cells.Item($rowA, $c).Value = $v
        $c++
    }
    $c = $firstCol
    foreach ($v in $valsA) {
        $ws.Cells.Item($rowB, $c).Value = $v
        $c++
    }
}

# Columns F..V are column indexes 6..22
$firstCol = 6
$lastCol = 22

Swap-RowData $ws 3 4 $firstCol $lastCol
Swap-RowData $ws 63 64 $firstCol $lastCol
Swap-RowData $ws 67 69 $firstCol $lastCol
Swap-RowData $ws 68 70 $firstCol $lastCol

function Add-MatchRow {
    param($ws, $templateRow, $newRow, $vals)

    # Clone formatting (bold/bordered index column, date-formatted column E, …)
    # from an existing fully-formatted data row.
    $ws.Range("A" + $templateRow + ":V" + $templateRow).Copy()
    $ws.Range("A" + $newRow + ":V" + $newRow).PasteSpecial(-4122)

    $c = 1
    foreach ($v in $vals) {
        $ws.Cells.Item($newRow, $c).Value = $v
        $c++
    }
}

Add-MatchRow $ws 79 80 @(79, "denmark", "2nd-division", "2023-2024", 45234.58333333334, "Brabrand", 2, "Nykobing", 1, 3.27, "03/11/2023 02:12", 3.29, "04/11/2023 13:47", 3.38, "03/11/2023 02:12", 3.54, "04/11/2023 13:47", 1.98, "03/11/2023 02:12", 2.06, "04/11/2023 13:47", "https://www.betexplorer.com/football/denmark/2nd-division/brabrand-nykobing/vkd3fa3o/")

Add-MatchRow $ws 79 81 @(80, "denmark", "2nd-division", "2023-2024", 45234.58333333334, "Aarhus Fremad", 3, "Hellerup", 2, 1.35, "03/11/2023 02:12", 1.33, "04/11/2023 13:39", 4.96, "03/11/2023 02:12", 5.5, "04/11/2023 13:39", 5.8, "03/11/2023 02:12", 7.43, "04/11/2023 13:39", "https://www.betexplorer.com/football/denmark/2nd-division/aarhus-fremad-hellerup/hh07gJIi/")

Add-MatchRow $ws 79 82 @(81, "denmark", "2nd-division", "2023-2024", 45234.58333333334, "Middelfart", 2, "FA 2000", 1, 1.5, "03/11/2023 02:12", 1.58, "04/11/2023 13:52", 4.13, "03/11/2023 02:12", 4.06, "04/11/2023 13:52", 5.07, "03/11/2023 02:12", 5.25, "04/11/2023 13:52", "https://www.betexplorer.com/football/denmark/2nd-division/middelfart-frederiksberg-alliancen-2000/nFaBhwYc/")

Add-MatchRow $ws 79 83 @(82, "denmark", "2nd-division", "2023-2024", 45234.79166666666, "AB Copenhagen", 1, "Roskilde", 1, 2.41, "03/11/2023 07:12", 2.79, "04/11/2023 18:59", 3.32, "03/11/2023 07:12", 3.68, "04/11/2023 18:59", 2.51, "03/11/2023 07:12", 2.27, "04/11/2023 18:59", "https://www.betexplorer.com/football/denmark/2nd-division/ab-copenhagen-roskilde/0KfimrdI/")
